# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# Change: cell B11 on sheet "Rules" goes from the text label "R40" to the
# text label "1" (a new, distinct shared string), while keeping the cell's
# existing formatting/style untouched.
#
# A plain  Range.Value = "1"  assignment would make Excel auto-detect the
# numeric-looking text and store it as a *number*, and pre-formatting the
# cell as Text (NumberFormat = "@") to force string storage would instead
# stamp the cell with a brand-new style. Neither matches the original
# cell's look (still text-typed, same style as before).
#
# So we stage the new text value in a scratch cell (forcing it to be text
# with a leading apostrophe), copy it, and use Paste Special > Values only
# into B11 - that brings over the *value + type* (text) without touching
# B11's existing number format/style. The scratch cell is then fully
# cleared (contents + formatting) so it leaves no trace on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$scratch = $ws.Range("Z1")
$scratch.Value = "'1"
$scratch.Copy()

$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues

$scratch.Clear()
